$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'39.585.76"
$ws.Range("E2").Value = "  -2.63%  "

# Row 3
$ws.Range("D3").Value = "'2.226.96"
$ws.Range("E3").Value = "  -6.06%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'297.76"
$ws.Range("E5").Value = "  -4.30%  "

# Row 6
$ws.Range("D6").Value = "'83.42"
$ws.Range("E6").Value = "  -3.51%  "

# Row 7
$ws.Range("D7").Value = "'0.513"
$ws.Range("E7").Value = "  -3.25%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "'0.471"
$ws.Range("E9").Value = "  -3.96%  "

# Row 10
$ws.Range("D10").Value = "'0.0777"
$ws.Range("E10").Value = "  -7.20%  "

# Row 11
$ws.Range("D11").Value = "'29.38"
$ws.Range("E11").Value = "  -2.93%  "

# Row 12
$ws.Range("D12").Value = "'47.76"
$ws.Range("E12").Value = "  -9.47%  "

# Row 13
$ws.Range("E13").Value = "  -2.06%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.34"
$ws.Range("E14").Value = "  -2.80%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'2.561.37"
$ws.Range("E15").Value = "  -6.28%  "

# Row 16
$ws.Range("D16").Value = "'14.20"
$ws.Range("E16").Value = "  -4.96%  "

# Row 17
$ws.Range("D17").Value = "'2.215.43"
$ws.Range("E17").Value = "  -6.19%  "

# Row 18
$ws.Range("D18").Value = "'0.720"
$ws.Range("E18").Value = "  -4.83%  "

# Row 19
$ws.Range("D19").Value = "'39.450.57"
$ws.Range("E19").Value = "  -2.76%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0878"
$ws.Range("E20").Value = "  -3.40%  "

# Row 21
$ws.Range("D21").Value = "'5.75"
$ws.Range("E21").Value = "  -6.02%  "

# Row 22
$ws.Range("D22").Value = "'65.24"
$ws.Range("E22").Value = "  -4.48%  "

# Row 23
$ws.Range("D23").Value = "'10.35"
$ws.Range("E23").Value = "  -3.52%  "

# Row 24
$ws.Range("D24").Value = "'228.49"
$ws.Range("E24").Value = "  -2.81%  "

# Row 25
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = "  -6.11%  "

# Row 27
$ws.Range("E27").Value = "  +1.07%  "

# Row 28
$ws.Range("D28").Value = "'22.71"
$ws.Range("E28").Value = "  -3.96%  "

# Row 29
$ws.Range("E29").Value = "  +0.44%  "

# Row 30
$ws.Range("D30").Value = "'9.15"
$ws.Range("E30").Value = "  -0.62%  "

# Row 31
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'32.20"
$ws.Range("E31").Value = "  -5.66%  "

# Row 32
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'149.93"
$ws.Range("E32").Value = "  -2.31%  "

# Row 33
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  -0.27%  "

# Row 34
$ws.Range("D34").Value = "'4.88"
$ws.Range("E34").Value = "  -5.72%  "

# Row 35
$ws.Range("D35").Value = "'0.0698"
$ws.Range("E35").Value = "  -4.03%  "

# Row 36
$ws.Range("E36").Value = "  -3.23%  "

# Row 37
$ws.Range("E37").Value = "  -3.32%  "

# Row 38
$ws.Range("D38").Value = "'0.0970"
$ws.Range("E38").Value = "  -2.97%  "

# Row 39
$ws.Range("D39").Value = "'15.33"
$ws.Range("E39").Value = "  -3.65%  "

# Row 40
$ws.Range("D40").Value = "'2.66"
$ws.Range("E40").Value = "  -3.90%  "

# Row 41
$ws.Range("D41").Value = "'1.66"
$ws.Range("E41").Value = "  -2.31%  "

# Row 42
$ws.Range("E42").Value = "  -4.05%  "

# Row 43
$ws.Range("D43").Value = "'1.914.13"
$ws.Range("E43").Value = "  -2.22%  "

# Row 44
$ws.Range("D44").Value = "'0.0260"
$ws.Range("E44").Value = "  -3.01%  "

# Row 45
$ws.Range("E45").Value = "  -15.23%  "

# Row 46
$ws.Range("D46").Value = "'16.25"
$ws.Range("E46").Value = "  -8.46%  "

# Row 47
$ws.Range("D47").Value = "'9.01"
$ws.Range("E47").Value = "  -3.43%  "

# Row 48
$ws.Range("D48").Value = "'2.63"
$ws.Range("E48").Value = "  -1.78%  "

# Row 49
$ws.Range("D49").Value = "'2.430.21"
$ws.Range("E49").Value = "  -6.53%  "

# Row 50
$ws.Range("D50").Value = "'70.97"
$ws.Range("E50").Value = "  -0.96%  "

# Row 51
$ws.Range("D51").Value = "'87.52"
$ws.Range("E51").Value = "  -5.86%  "
